# Updated symbol list on Wed Dec 28 22:54:22 UTC 2022 with GitHub Actions
#
# Applies updated crypto price/volume data to the sheet. Numeric-looking
# values are stored as text (matching the original inlineStr cells), so
# they are written with a leading apostrophe to force Excel to keep them
# as strings (preserving formatting such as leading/trailing zeros)
# instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple numeric-looking price/volume updates (row unchanged)
$ws.Range("D2").Value  = "'243.27"
$ws.Range("D3").Value  = "'23.70"
$ws.Range("D4").Value  = "'5.232"
$ws.Range("D7").Value  = "'3.229"
$ws.Range("D9").Value  = "'0.8894"
$ws.Range("D11").Value = "'0.07076"
$ws.Range("D14").Value = "'0.09311"
$ws.Range("D15").Value = "'3.814"
$ws.Range("D16").Value = "'0.001517"
$ws.Range("D17").Value = "'0.04717"

# Rows 18-24: coin ranking list shifted by one position, with new price data
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006144"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").Value = "'0.001248"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004061"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.00008695"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.545"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.133"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01011"
$ws.Range("E24").Value = "23OneONEBestin24h"

# Remaining simple numeric-looking price/volume updates
$ws.Range("D26").Value = "'0.1309"
$ws.Range("D28").Value = "'0.0002327"
$ws.Range("D40").Value = "'0.03716"
$ws.Range("D41").Value = "'0.006246"
$ws.Range("D42").Value = "'0.1043"
$ws.Range("D43").Value = "'0.002500"
$ws.Range("D44").Value = "'0.007136"
$ws.Range("D45").Value = "'0.00005323"

$ws.Range("D47").Value = "'0.5347"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

$ws.Range("D48").Value = "'0.002323"
$ws.Range("E48").Value = "47BOLOBOLO"

$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.0001999"
